# Updates cryptos list values (price + 1h volume change) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values are plain decimal numbers (e.g. "305.28") must be
# forced to Text so Excel keeps them as literal strings (matching the existing column
# formatting of mixed "42.667.19"-style grouped values) instead of auto-converting them
# to numeric values.
$textForcedCells = @("D5", "D6", "D10", "D12", "D14", "D19", "D21", "D22", "D23", "D25", "D26", "D28", "D29", "D30", "D34", "D36", "D45", "D46", "D47", "D49", "D50")
foreach ($cellRef in $textForcedCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "42.655.24"
$ws.Range("D3").Value = "2.279.54"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "305.28"
$ws.Range("E5").Value = "  +1.74%  "
$ws.Range("D6").Value = "96.57"
$ws.Range("E6").Value = "  -1.27%  "
$ws.Range("E7").Value = "  -2.75%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -3.17%  "
$ws.Range("D10").Value = "35.57"
$ws.Range("E10").Value = "  -2.06%  "
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "18.26"
$ws.Range("E12").Value = "  +2.66%  "
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").Value = "6.75"
$ws.Range("E14").Value = "  -2.10%  "
$ws.Range("D15").Value = "2.633.98"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").Value = "2.292.16"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "42.589.14"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").Value = "12.90"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("D21").Value = "6.01"
$ws.Range("E21").Value = "  -2.01%  "
$ws.Range("D22").Value = "67.08"
$ws.Range("E22").Value = "  -1.76%  "
$ws.Range("D23").Value = "236.17"
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("E24").Value = "  -2.93%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "2.46"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "25.11"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").Value = "165.88"
$ws.Range("E29").Value = "  +1.68%  "
$ws.Range("D30").Value = "2.06"
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("E31").Value = "  -1.29%  "
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").Value = "4.76"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("E35").Value = "  -2.95%  "
$ws.Range("D36").Value = "17.60"
$ws.Range("E36").Value = "  -3.25%  "
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("E40").Value = "  -2.29%  "
$ws.Range("E41").Value = "  -1.66%  "
$ws.Range("E42").Value = "  -2.99%  "
$ws.Range("D43").Value = "2.001.12"
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("E44").Value = "  -2.54%  "
$ws.Range("D45").Value = "18.13"
$ws.Range("E45").Value = "  +3.39%  "
$ws.Range("D46").Value = "9.97"
$ws.Range("E46").Value = "  -3.36%  "
$ws.Range("D47").Value = "2.09"
$ws.Range("E47").Value = "  -6.73%  "
$ws.Range("E48").Value = "  -2.58%  "
$ws.Range("D49").Value = "2.85"
$ws.Range("E49").Value = "  +4.56%  "
$ws.Range("D50").Value = "53.53"
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("D51").Value = "2.502.48"
$ws.Range("E51").Value = "  -1.16%  "
